$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving number format on cells whose target value is a numeric-looking string
$textCells = @("G9", "G10", "G11", "G13", "G14", "G15", "D16", "G16", "D17", "G17", "G18", "G19", "G21", "G22", "D23", "G23", "D24", "G24", "D25", "G25", "D26", "G26", "D27", "D28", "G28", "G30", "D31", "D32", "G32", "D33", "G33", "D34", "D35", "G35", "D36", "G36", "D37", "D38", "D39", "G41", "H41", "G43", "H43")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Write values
$ws.Range("C8").Value = 38
$ws.Range("C9").Value = 97
$ws.Range("G9").Value = '24832.00'
$ws.Range("C10").Value = 41
$ws.Range("G10").Value = '19352.00'
$ws.Range("C11").Value = 85
$ws.Range("G11").Value = '56270.00'
$ws.Range("C12").Value = 61
$ws.Range("C13").Value = 58
$ws.Range("G13").Value = '7888.00'
$ws.Range("C14").Value = 17
$ws.Range("G14").Value = '391.00'
$ws.Range("C15").Value = 26
$ws.Range("G15").Value = '1300.00'
$ws.Range("C16").Value = 74
$ws.Range("D16").Value = '5.0'
$ws.Range("E16").Value = 'Providing & Fixing of  of 3/5 pin 6 amp. flush type  non modular socket  made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = '2442.00'
$ws.Range("C17").Value = 30
$ws.Range("D17").Value = '6.0'
$ws.Range("E17").Value = 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F17").Value = 78
$ws.Range("G17").Value = '2340.00'
$ws.Range("C18").Value = 21
$ws.Range("G18").Value = '4599.00'
$ws.Range("C19").Value = 35
$ws.Range("G19").Value = '10605.00'
$ws.Range("C20").Value = 94
$ws.Range("C21").Value = 50
$ws.Range("G21").Value = '2000.00'
$ws.Range("C22").Value = 55
$ws.Range("G22").Value = '3080.00'
$ws.Range("A23").Value = 'Mtr.'
$ws.Range("C23").Value = 93
$ws.Range("D23").Value = '20'
$ws.Range("E23").Value = '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
$ws.Range("F23").Value = 122
$ws.Range("G23").Value = '11346.00'
$ws.Range("A24").Value = 'Set'
$ws.Range("C24").Value = 93
$ws.Range("D24").Value = '13.0'
$ws.Range("E24").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F24").Value = 5733
$ws.Range("G24").Value = '533169.00'
$ws.Range("A25").Value = ''
$ws.Range("C25").Value = 60
$ws.Range("D25").Value = '14.0'
$ws.Range("E25").Value = 'Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = '0.00'
$ws.Range("A26").Value = 'Mtr.'
$ws.Range("C26").Value = 89
$ws.Range("D26").Value = '23'
$ws.Range("E26").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = '1780.00'
$ws.Range("C27").Value = 87
$ws.Range("D27").Value = '15.0'
$ws.Range("E27").Value = 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("A28").Value = 'Each'
$ws.Range("C28").Value = 54
$ws.Range("D28").Value = '25'
$ws.Range("E28").Value = '1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )'
$ws.Range("F28").Value = 1890
$ws.Range("G28").Value = '102060.00'
$ws.Range("C29").Value = 55
$ws.Range("C30").Value = 94
$ws.Range("G30").Value = '46248.00'
$ws.Range("C31").Value = 99
$ws.Range("D31").Value = '29'
$ws.Range("E31").Value = 'Single pole MCB   (With B/C curve tripping Characteristics)'
$ws.Range("A32").Value = 'Each'
$ws.Range("C32").Value = 64
$ws.Range("D32").Value = '30'
$ws.Range("E32").Value = ' 6 A to 32 A rating'
$ws.Range("F32").Value = 187
$ws.Range("G32").Value = '11968.00'
$ws.Range("C33").Value = 37
$ws.Range("D33").Value = '32'
$ws.Range("E33").Value = ' 50/63 A rating'
$ws.Range("F33").Value = 900
$ws.Range("G33").Value = '33300.00'
$ws.Range("C34").Value = 23
$ws.Range("D34").Value = '18.0'
$ws.Range("E34").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("A35").Value = ''
$ws.Range("C35").Value = 69
$ws.Range("D35").Value = '34'
$ws.Range("E35").Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = '0.00'
$ws.Range("A36").Value = 'Each'
$ws.Range("C36").Value = 77
$ws.Range("D36").Value = '35'
$ws.Range("E36").Value = '8 Way (8+2)'
$ws.Range("F36").Value = 2184
$ws.Range("G36").Value = '168168.00'
$ws.Range("C37").Value = 91
$ws.Range("D37").Value = '36'
$ws.Range("E37").Value = 'Total'
$ws.Range("A38").Value = '%'
$ws.Range("C38").Value = 31
$ws.Range("D38").Value = '37'
$ws.Range("E38").Value = 'Add Tender Premium '
$ws.Range("A39").Value = ''
$ws.Range("C39").Value = 36
$ws.Range("D39").Value = '38'
$ws.Range("E39").Value = 'Grand Total'
$ws.Range("G41").Value = '1043138.00'
$ws.Range("H41").Value = '1043138.00'
$ws.Range("G43").Value = '1043138.00'
$ws.Range("H43").Value = '1043138.00'

# Restore default (General) number format on the forced cells so formatting matches the rest of the sheet
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
}
